$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# C2: explanatory note (same "column C note" pattern as C1)
# -----------------------------------------------------------------
$ws.Range("C2").Value = "The last 3 categories are all from me from combining the original categories"

# -----------------------------------------------------------------
# New rows 21-23: the three combined/derived cause categories
# -----------------------------------------------------------------
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "Manmade"
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "Natural"
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "Miscellaneous/Unknown"

# Match the look of the rest of the table (font/fill/border) by copying
# the format of an existing data row, then trim the border down to just
# left/right (no top/bottom) to match the "inner" rows of this new block.
$ws.Range("A17:B17").Copy() | Out-Null
$ws.Range("A21:B21").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A22:B22").PasteSpecial(-4122) | Out-Null
$ws.Range("A23:B23").PasteSpecial(-4122) | Out-Null

foreach ($addr in @("A21","B21","A22","B22","A23","B23")) {
    $cell = $ws.Range($addr)
    $cell.Borders.Item(8).LineStyle = -4142   # xlEdgeTop    -> none
    $cell.Borders.Item(9).LineStyle = -4142   # xlEdgeBottom -> none
}

# -----------------------------------------------------------------
# Row heights: rows 19/20 shrink (re-wrap), new row 23 wraps to 2 lines
# -----------------------------------------------------------------
$ws.Rows.Item(19).RowHeight = 28.8
$ws.Rows.Item(20).RowHeight = 28.8
$ws.Rows.Item(23).RowHeight = 28.8

# -----------------------------------------------------------------
# View state: selection moved to B23, scrolled down toward the bottom
# -----------------------------------------------------------------
$wnd = $excel.ActiveWindow
$wnd.ScrollRow = 10
$wnd.ScrollColumn = 1
$ws.Range("B23").Select() | Out-Null

Write-Host "Edit complete"
